$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4 (Jalen Brunson, Desmond Bane, Devin Vassell) are unaffected.
# Starting at row 5 a new player (Norman Powell) is inserted and the rest
# of the roster shifts down / gets reshuffled, so rewrite columns A
# (Oyuncu Adi / player), B (Pozisyon / position) and C (Takim / team)
# for rows 5-19 with their final values.

$ws.Range("A5").Value = "Norman Powell"
$ws.Range("B5").Value = "SG,SF"
$ws.Range("C5").Value = "LA Clippers"

$ws.Range("A6").Value = "Shaedon Sharpe"
$ws.Range("B6").Value = "SG,SF"
$ws.Range("C6").Value = "Portland Trail Blazers"

$ws.Range("A7").Value = "LeBron James"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Los Angeles Lakers"

$ws.Range("A8").Value = "Dalton Knecht"
$ws.Range("B8").Value = "SG,SF"
$ws.Range("C8").Value = "Los Angeles Lakers"

$ws.Range("A9").Value = "P.J. Washington"
$ws.Range("B9").Value = "PF"
$ws.Range("C9").Value = "Dallas Mavericks"

$ws.Range("A10").Value = "Walker Kessler"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Utah Jazz"

$ws.Range("A11").Value = "Alperen Sengün"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Houston Rockets"

$ws.Range("A12").Value = "Nicolas Claxton"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "Brooklyn Nets"

$ws.Range("A13").Value = "Moussa Diabate"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Charlotte Hornets"

$ws.Range("A14").Value = "Trae Young"
$ws.Range("B14").Value = "PG"
$ws.Range("C14").Value = "Atlanta Hawks"

$ws.Range("A15").Value = "Devin Booker"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Phoenix Suns"

$ws.Range("A16").Value = "Coby White"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Chicago Bulls"

$ws.Range("A17").Value = "Immanuel Quickley"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("C17").Value = "Toronto Raptors"

$ws.Range("A18").Value = "Kawhi Leonard"
$ws.Range("B18").Value = "SG,SF,PF"
$ws.Range("C18").Value = "LA Clippers"

$ws.Range("A19").Value = "Dereck Lively II"
$ws.Range("B19").Value = "C"
$ws.Range("C19").Value = "Dallas Mavericks"
